$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new "Price" values below are plain decimal numbers (e.g. "232.46").
# If written to a General-formatted cell, Excel auto-converts them to a
# Number and silently drops significant trailing zeros (e.g. "56.70" -> 56.7).
# The source data keeps these as text, so force Text format first.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D16", "D17", "D21", "D22", "D23", "D29", "D30", "D31", "D34", "D35", "D38", "D41", "D43", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.363.81'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '2.332.67'
$ws.Range('E3').Value = '  +2.66%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = '232.46'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').Value = '0.648'
$ws.Range('E6').Value = '  +2.01%  '
$ws.Range('D7').Value = '66.54'
$ws.Range('E7').Value = '  +3.79%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '0.454'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '0.0968'
$ws.Range('E10').Value = '  -3.65%  '
$ws.Range('D11').Value = '56.70'
$ws.Range('E11').Value = '  -1.19%  '
$ws.Range('D12').Value = '26.81'
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').Value = '2.678.62'
$ws.Range('E13').Value = '  +2.53%  '
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').Value = '6.26'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').Value = '0.854'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '2.318.68'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').Value = '43.261.48'
$ws.Range('E19').Value = '  -1.65%  '
$ws.Range('E20').Value = '  -2.90%  '
$ws.Range('D21').Value = '74.18'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = '6.24'
$ws.Range('E22').Value = '  +2.12%  '
$ws.Range('D23').Value = '248.99'
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('E24').Value = '  +13.35%  '
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '22.21'
$ws.Range('E29').Value = '  +6.60%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '174.17'
$ws.Range('E30').Value = '  +1.60%  '
$ws.Range('D31').Value = '1.47'
$ws.Range('E31').Value = '  +5.74%  '
$ws.Range('E32').Value = '  -7.61%  '
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('D34').Value = '5.02'
$ws.Range('E34').Value = '  +4.10%  '
$ws.Range('D35').Value = '0.0690'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('E37').Value = '  +8.64%  '
$ws.Range('D38').Value = '6.50'
$ws.Range('E38').Value = '  -0.82%  '
$ws.Range('E39').Value = '  -4.77%  '
$ws.Range('E40').Value = '  -2.36%  '
$ws.Range('D41').Value = '9.10'
$ws.Range('E41').Value = '  +10.14%  '
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '18.37'
$ws.Range('E43').Value = '  +4.25%  '
$ws.Range('E44').Value = '  +8.22%  '
$ws.Range('D45').Value = '99.44'
$ws.Range('E45').Value = '  +1.11%  '
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').Value = '0.0946'
$ws.Range('E47').Value = '  -3.49%  '
$ws.Range('D48').Value = '4.34'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('D49').Value = '1.446.95'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = '9.96'
$ws.Range('E50').Value = '  -5.48%  '
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D51').Value = '0.000205'
$ws.Range('E51').Value = '  -13.58%  '
